$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 804 ("Hortaliza, Terminal Hortofrutícola
# Agro Chillán - Betarraga" weekly price update). This shifts the former rows
# 804:822 down to 806:824.
$ws.Rows("804:805").Insert()

# New row 804
$ws.Cells.Item(804, 1).Value = 7
$ws.Cells.Item(804, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(804, 3).Value = "Ñuble"
$ws.Cells.Item(804, 4).Value = 45239
$ws.Cells.Item(804, 5).Value = 16
$ws.Cells.Item(804, 6).Value = 100114014
$ws.Cells.Item(804, 7).Value = "Betarraga"
$ws.Cells.Item(804, 8).Value = "Sin especificar"
$ws.Cells.Item(804, 9).Value = "Primera"
$ws.Cells.Item(804, 10).Value = 300
$ws.Cells.Item(804, 11).Value = 900
$ws.Cells.Item(804, 12).Value = 900
$ws.Cells.Item(804, 13).Value = 900
$ws.Cells.Item(804, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(804, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(804, 16).Value = 180
$ws.Cells.Item(804, 17).Value = 5
$ws.Cells.Item(804, 18).Value = "Hortaliza"

# New row 805
$ws.Cells.Item(805, 1).Value = 7
$ws.Cells.Item(805, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(805, 3).Value = "Ñuble"
$ws.Cells.Item(805, 4).Value = 45239
$ws.Cells.Item(805, 5).Value = 16
$ws.Cells.Item(805, 6).Value = 100114014
$ws.Cells.Item(805, 7).Value = "Betarraga"
$ws.Cells.Item(805, 8).Value = "Sin especificar"
$ws.Cells.Item(805, 9).Value = "Primera"
$ws.Cells.Item(805, 10).Value = 500
$ws.Cells.Item(805, 11).Value = 700
$ws.Cells.Item(805, 12).Value = 700
$ws.Cells.Item(805, 13).Value = 700
$ws.Cells.Item(805, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(805, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(805, 16).Value = 140
$ws.Cells.Item(805, 17).Value = 5
$ws.Cells.Item(805, 18).Value = "Hortaliza"

# Make sure the date cells keep the date/time number format used elsewhere in
# column D (style index 2 in the original workbook).
$ws.Range("D804:D805").NumberFormat = "YYYY-MM-DD HH:MM:SS"
